$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The header row held template placeholders for the "is_locked" and
# "is_enabled" dictionary columns in D1:E1 (with the previous F1/G1 holding
# "order_by"/"rem"). Those two editable-lookup columns are no longer part of
# the import template, so delete D1:E1 and shift the remaining cells
# (order_by, rem) left into their place.
$ws.Range("D1:E1").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
